$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The deliverables list needs a new row ("Progress Report") inserted right
# before the existing "Handover Document" row (row 17), pushing that row
# (and its "No." numbering) down to row 18.

# Insert a blank row at row 17, shifting "Handover Document" (and below) down.
$ws.Rows.Item(17).Insert()

# Copy formatting from the row above (row 16) onto the newly inserted row so
# it keeps the same borders/style as the rest of the table instead of
# Excel's default blank-row formatting.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "Progress Report" row.
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Progress Report"
$ws.Cells.Item(17, 3).Value = "Yes"

# Renumber the "Handover Document" row that got pushed down to row 18.
$ws.Cells.Item(18, 1).Value = 17

# Match the author's final selection.
$ws.Range("C17").Select()
